# Update gh-pages to output generated at 456a3b4
# Applies numeric "想去人数" (F column) bumps and refreshed "Cover" (I column)
# image URLs to both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# F-column (numeric) updates: row -> new value
$fUpdates = @{
    2  = 40
    3  = 21
    7  = 106
    8  = 76
    9  = 445
    10 = 43
    12 = 548
    14 = 284
    16 = 346
    17 = 112
    21 = 87
    22 = 857
    23 = 1377
    25 = 309
    33 = 260
    34 = 1595
    38 = 567
    40 = 3430
    41 = 406
    42 = 179
    43 = 879
    45 = 58
}

# I-column (Cover URL) updates: row -> new value
$iUpdates = @{
    3 = "//i2.hdslb.com/bfs/openplatform/202401/8YmblqtV1706524591857.png"
    8 = "//i0.hdslb.com/bfs/openplatform/202401/BjtgGUbI1706525642100.png"
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    foreach ($row in $iUpdates.Keys) {
        $ws.Range("I$row").Value = $iUpdates[$row]
    }
}
